$d = $word.ActiveDocument

$replacements = @(
    @("29×66=", "75×40="),
    @("20×16=", "98×19="),
    @("48×66=", "26×69="),
    @("65×99=", "25×43="),
    @("20×97=", "31×90="),
    @("24×56=", "83×22="),
    @("64×50=", "20×82="),
    @("31×78=", "16×95="),
    @("69×75=", "43×82="),
    @("81×62=", "42×53="),
    @("71×34=", "15×72="),
    @("43×12=", "74×51="),
    @("87×86=", "94×42="),
    @("22×27=", "24×79="),
    @("79×45=", "65×74="),
    @("13×75=", "46×90="),
    @("70×68=", "21×92="),
    @("58×43=", "33×56="),
    @("78×19=", "79×62="),
    @("86×11=", "41×33="),
    @("32×39=", "42×19="),
    @("77×74=", "31×41="),
    @("58×51=", "97×79="),
    @("15×35=", "48×59="),
    @("44×57=", "33×75=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
